$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-08-27 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-08-28 Thursday", 2)

$d.Content.Find.Execute("87×54=4698", $true, $false, $false, $false, $false, $true, 1, $false, "37×38=1406", 2)
$d.Content.Find.Execute("68×71=4828", $true, $false, $false, $false, $false, $true, 1, $false, "11×13=143", 2)
$d.Content.Find.Execute("79×35=2765", $true, $false, $false, $false, $false, $true, 1, $false, "33×45=1485", 2)
$d.Content.Find.Execute("12×27=324", $true, $false, $false, $false, $false, $true, 1, $false, "49×86=4214", 2)
$d.Content.Find.Execute("82×18=1476", $true, $false, $false, $false, $false, $true, 1, $false, "86×49=4214", 2)

$d.Content.Find.Execute("55×50=2750", $true, $false, $false, $false, $false, $true, 1, $false, "61×28=1708", 2)
$d.Content.Find.Execute("86×62=5332", $true, $false, $false, $false, $false, $true, 1, $false, "42×82=3444", 2)
$d.Content.Find.Execute("35×23=805", $true, $false, $false, $false, $false, $true, 1, $false, "74×82=6068", 2)
$d.Content.Find.Execute("25×74=1850", $true, $false, $false, $false, $false, $true, 1, $false, "93×75=6975", 2)
$d.Content.Find.Execute("36×52=1872", $true, $false, $false, $false, $false, $true, 1, $false, "96×34=3264", 2)

$d.Content.Find.Execute("30×37=1110", $true, $false, $false, $false, $false, $true, 1, $false, "93×22=2046", 2)
$d.Content.Find.Execute("21×40=840", $true, $false, $false, $false, $false, $true, 1, $false, "89×12=1068", 2)
$d.Content.Find.Execute("30×23=690", $true, $false, $false, $false, $false, $true, 1, $false, "40×90=3600", 2)
$d.Content.Find.Execute("90×93=8370", $true, $false, $false, $false, $false, $true, 1, $false, "85×31=2635", 2)
$d.Content.Find.Execute("98×11=1078", $true, $false, $false, $false, $false, $true, 1, $false, "14×45=630", 2)

$d.Content.Find.Execute("53×82=4346", $true, $false, $false, $false, $false, $true, 1, $false, "18×87=1566", 2)
$d.Content.Find.Execute("18×38=684", $true, $false, $false, $false, $false, $true, 1, $false, "83×13=1079", 2)
$d.Content.Find.Execute("83×96=7968", $true, $false, $false, $false, $false, $true, 1, $false, "91×37=3367", 2)
$d.Content.Find.Execute("90×96=8640", $true, $false, $false, $false, $false, $true, 1, $false, "98×25=2450", 2)
$d.Content.Find.Execute("18×29=522", $true, $false, $false, $false, $false, $true, 1, $false, "81×63=5103", 2)

$d.Content.Find.Execute("79×37=2923", $true, $false, $false, $false, $false, $true, 1, $false, "28×76=2128", 2)
$d.Content.Find.Execute("15×60=900", $true, $false, $false, $false, $false, $true, 1, $false, "22×59=1298", 2)
$d.Content.Find.Execute("20×46=920", $true, $false, $false, $false, $false, $true, 1, $false, "82×45=3690", 2)
$d.Content.Find.Execute("93×51=4743", $true, $false, $false, $false, $false, $true, 1, $false, "76×26=1976", 2)
$d.Content.Find.Execute("92×20=1840", $true, $false, $false, $false, $false, $true, 1, $false, "46×31=1426", 2)
